# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G (header "K") values for rows 2-6 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 3
